$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Title cells (C3 / J3) become bold ---
$ws.Range("C3").Font.Bold = $true
$ws.Range("J3").Font.Bold = $true

# --- Update the "value" cells (both the left table and its mirrored
#     right-hand copy) from their old text to the new "asd" text ---
$ws.Range("C7").Value = "asd"
$ws.Range("C8").Value = "asd"
$ws.Range("C9").Value = "asd"
$ws.Range("J7").Value = "asd"
$ws.Range("J8").Value = "asd"
$ws.Range("J9").Value = "asd"

# --- Merge the label (A:B) and value (C:F) blocks on rows 7-9 ---
$ws.Range("A7:B7").Merge()
$ws.Range("A8:B8").Merge()
$ws.Range("A9:B9").Merge()
$ws.Range("C7:F7").Merge()
$ws.Range("C8:F8").Merge()
$ws.Range("C9:F9").Merge()

# --- Apply a thin box border + left alignment around those merged
#     blocks and around the matching H column cells ---
$labels = $ws.Range("A7:B9")
$labels.Borders.LineStyle = 1
$labels.HorizontalAlignment = -4131

$values = $ws.Range("C7:F9")
$values.Borders.LineStyle = 1
$values.HorizontalAlignment = -4131

$hcol = $ws.Range("H7:H9")
$hcol.Borders.LineStyle = 1
$hcol.HorizontalAlignment = -4131
